$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several rows had "Primary topic" (col B) and "Secondary topic" (col C)
# swapped by mistake -- "Leadership content" was recorded in the wrong
# column. Fix the golden data by swapping B/C for the affected rows.
$rows = 7, 10, 11, 30, 33, 41, 49, 55, 57

foreach ($r in $rows) {
    $bCell = $ws.Cells.Item($r, 2)
    $cCell = $ws.Cells.Item($r, 3)

    $bVal = $bCell.Value()
    $cVal = $cCell.Value()

    $bCell.Value = $cVal
    $cCell.Value = $bVal
}

# Rows 33 and 55 had the (red) highlight font applied to the "Leadership
# content" cell in column C. After the swap, column B now holds
# "Leadership content" and should no longer carry that highlight -- the
# highlight formatting does not travel with the swapped value.
$ws.Cells.Item(33, 2).ClearFormats()
$ws.Cells.Item(55, 3).ClearFormats()

# Move the active selection from A30 to B30.
$ws.Range("B30").Select()
